$wb = $excel.ActiveWorkbook

$wsALC = $wb.Worksheets.Item("ALC")
$wsARM = $wb.Worksheets.Item("ARM")
$wsBSM = $wb.Worksheets.Item("BSM")
$wsCRP = $wb.Worksheets.Item("CRP")
$wsCUL = $wb.Worksheets.Item("CUL")
$wsGSM = $wb.Worksheets.Item("GSM")
$wsLTW = $wb.Worksheets.Item("LTW")
$wsWVR = $wb.Worksheets.Item("WVR")

# ALC!row88
$wsALC.Range("H88").Value = 4170058.8
$wsALC.Range("I88").Value = 14287652
$wsALC.Range("J88").Value = 3990.8823
$wsALC.Range("K88").Value = 14287652
$wsALC.Range("L88").Value = 3990.8823
$wsALC.Range("M88").Value = -14287246
$wsALC.Range("N88").Value = -4802.8823

# ALC!row91
$wsALC.Range("H91").Value = 4170058.8
$wsALC.Range("I91").Value = 14287652
$wsALC.Range("J91").Value = 3990.8823
$wsALC.Range("K91").Value = 14287652
$wsALC.Range("L91").Value = 3990.8823
$wsALC.Range("M91").Value = -14286248
$wsALC.Range("N91").Value = -6798.8823

# ALC!row113
$wsALC.Range("H113").Value = 3665.3333
$wsALC.Range("I113").Value = 3330.6667
$wsALC.Range("K113").Value = 3330.6667
$wsALC.Range("M113").Value = -76.66670000000022

# ALC!row132
$wsALC.Range("H132").Value = 5612.6787
$wsALC.Range("I132").Value = 6379.5415
$wsALC.Range("J132").Value = 1011.5
$wsALC.Range("K132").Value = 19138.6245
$wsALC.Range("L132").Value = 3034.5
$wsALC.Range("M132").Value = -16608.6245
$wsALC.Range("N132").Value = -8094.5

# ARM!row32
$wsARM.Range("H32").Value = 1471.8857
$wsARM.Range("I32").Value = 1620.7742
$wsARM.Range("K32").Value = 1620.7742
$wsARM.Range("M32").Value = -1333.7742

# ARM!row45
$wsARM.Range("H45").Value = 43924.6
$wsARM.Range("I45").Value = 48535.11
$wsARM.Range("K45").Value = 48535.11
$wsARM.Range("M45").Value = -48158.11

# ARM!row74
$wsARM.Range("H74").Value = 207290.19
$wsARM.Range("I74").Value = 265794.8
$wsARM.Range("J74").Value = 2524
$wsARM.Range("K74").Value = 265794.8
$wsARM.Range("L74").Value = 2524
$wsARM.Range("M74").Value = -264920.8
$wsARM.Range("N74").Value = -4272

# ARM!row77
$wsARM.Range("H77").Value = 207290.19
$wsARM.Range("I77").Value = 265794.8
$wsARM.Range("J77").Value = 2524
$wsARM.Range("K77").Value = 1328974
$wsARM.Range("L77").Value = 12620
$wsARM.Range("M77").Value = -1324606
$wsARM.Range("N77").Value = -21356

# ARM!row122
$wsARM.Range("H122").Value = 4471.07
$wsARM.Range("I122").Value = 4942.0938
$wsARM.Range("J122").Value = 3100.818
$wsARM.Range("K122").Value = 14826.2814
$wsARM.Range("L122").Value = 9302.454000000002
$wsARM.Range("M122").Value = -12376.2814
$wsARM.Range("N122").Value = -14202.454

# ARM!row133
$wsARM.Range("H133").Value = 105998.5
$wsARM.Range("J133").Value = 105998.5
$wsARM.Range("L133").Value = 105998.5
$wsARM.Range("N133").Value = -111058.5

# BSM!row105
$wsBSM.Range("H105").Value = 23638562
$wsBSM.Range("I105").Value = 2501699.8
$wsBSM.Range("J105").Value = 35716770
$wsBSM.Range("K105").Value = 2501699.8
$wsBSM.Range("L105").Value = 35716770
$wsBSM.Range("M105").Value = -2499952.8
$wsBSM.Range("N105").Value = -35720264

# CRP!row31
$wsCRP.Range("H31").Value = 2721302
$wsCRP.Range("I31").Value = 3103.7856
$wsCRP.Range("J31").Value = 6949610.5
$wsCRP.Range("K31").Value = 3103.7856
$wsCRP.Range("L31").Value = 6949610.5
$wsCRP.Range("M31").Value = -2808.7856
$wsCRP.Range("N31").Value = -6950200.5

# CRP!row34
$wsCRP.Range("H34").Value = 2721302
$wsCRP.Range("I34").Value = 3103.7856
$wsCRP.Range("J34").Value = 6949610.5
$wsCRP.Range("K34").Value = 3103.7856
$wsCRP.Range("L34").Value = 6949610.5
$wsCRP.Range("M34").Value = -2901.7856
$wsCRP.Range("N34").Value = -6950014.5

# CRP!row58
$wsCRP.Range("H58").Value = 2463.9033
$wsCRP.Range("I58").Value = 1850.6316
$wsCRP.Range("K58").Value = 1850.6316
$wsCRP.Range("M58").Value = -1647.6316

# CRP!row99
$wsCRP.Range("H99").Value = 3039.125
$wsCRP.Range("I99").Value = 3414.3333
$wsCRP.Range("K99").Value = 3414.3333
$wsCRP.Range("M99").Value = -1916.3333

# CRP!row122
$wsCRP.Range("H122").Value = 2551.4375
$wsCRP.Range("I122").Value = 2381.4
$wsCRP.Range("K122").Value = 7144.200000000001
$wsCRP.Range("M122").Value = -4694.200000000001

# CRP!row126
$wsCRP.Range("H126").Value = 3039.125
$wsCRP.Range("I126").Value = 3414.3333
$wsCRP.Range("K126").Value = 10242.9999
$wsCRP.Range("M126").Value = -7772.999899999999

# CRP!row136
$wsCRP.Range("H136").Value = 2463.9033
$wsCRP.Range("I136").Value = 1850.6316
$wsCRP.Range("K136").Value = 5551.8948
$wsCRP.Range("M136").Value = -3001.8948

# CUL!row38
$wsCUL.Range("H38").Value = 193.07143
$wsCUL.Range("J38").Value = 260
$wsCUL.Range("L38").Value = 780
$wsCUL.Range("N38").Value = -1474

# CUL!row132
$wsCUL.Range("H132").Value = 1948.4615
$wsCUL.Range("I132").Value = 995
$wsCUL.Range("K132").Value = 8955
$wsCUL.Range("M132").Value = -6425

# CUL!row134
$wsCUL.Range("H134").Value = 2874.7856
$wsCUL.Range("I134").Value = 2025.4
$wsCUL.Range("K134").Value = 6076.200000000001
$wsCUL.Range("M134").Value = -1006.200000000001

# GSM!row7
$wsGSM.Range("H7").Value = 10000
$wsGSM.Range("J7").Value = 10000
$wsGSM.Range("L7").Value = 10000
$wsGSM.Range("N7").Value = -10224

# GSM!row8
$wsGSM.Range("H8").Value = 10000
$wsGSM.Range("J8").Value = 10000
$wsGSM.Range("L8").Value = 10000
$wsGSM.Range("N8").Value = -10278

# GSM!row96
$wsGSM.Range("H96").Value = 46824
$wsGSM.Range("J96").Value = 46824
$wsGSM.Range("L96").Value = 46824
$wsGSM.Range("N96").Value = -52316

# GSM!row102
$wsGSM.Range("H102").Value = 11197.667
$wsGSM.Range("I102").Value = 2785.8
$wsGSM.Range("J102").Value = 53257
$wsGSM.Range("K102").Value = 2785.8
$wsGSM.Range("L102").Value = 53257
$wsGSM.Range("M102").Value = -1163.8
$wsGSM.Range("N102").Value = -56501

# GSM!row122
$wsGSM.Range("H122").Value = 3038.1428
$wsGSM.Range("I122").Value = 2358.5386
$wsGSM.Range("J122").Value = 5001.4443
$wsGSM.Range("K122").Value = 7075.6158
$wsGSM.Range("L122").Value = 15004.3329
$wsGSM.Range("M122").Value = -4625.6158
$wsGSM.Range("N122").Value = -19904.3329

# LTW!row93
$wsLTW.Range("H93").Value = 2093.2856
$wsLTW.Range("I93").Value = 1755.091
$wsLTW.Range("K93").Value = 1755.091
$wsLTW.Range("M93").Value = -507.0909999999999

# LTW!row122
$wsLTW.Range("H122").Value = 7029.9165
$wsLTW.Range("I122").Value = 4040.111
$wsLTW.Range("J122").Value = 15999.333
$wsLTW.Range("K122").Value = 12120.333
$wsLTW.Range("L122").Value = 47997.999
$wsLTW.Range("M122").Value = -9670.332999999999
$wsLTW.Range("N122").Value = -52897.999

# WVR!row22
$wsWVR.Range("H22").Value = 9004.333000000001
$wsWVR.Range("I22").Value = 9004.333000000001
$wsWVR.Range("J22").Value = 0
$wsWVR.Range("K22").Value = 9004.333000000001
$wsWVR.Range("L22").Value = 0
$wsWVR.Range("N22").ClearContents()
$wsWVR.Range("M22").Value = -8711.333000000001

# WVR!row54
$wsWVR.Range("H54").Value = 46627.75
$wsWVR.Range("J54").Value = 48993.145
$wsWVR.Range("L54").Value = 48993.145
$wsWVR.Range("N54").Value = -50033.145

# WVR!row107
$wsWVR.Range("H107").Value = 671.90625
$wsWVR.Range("I107").Value = 570.619
$wsWVR.Range("J107").Value = 865.2727
$wsWVR.Range("K107").Value = 1711.857
$wsWVR.Range("L107").Value = 2595.8181
$wsWVR.Range("M107").Value = 208.143
$wsWVR.Range("N107").Value = -6435.8181

# WVR!row136
$wsWVR.Range("H136").Value = 197129.3
$wsWVR.Range("I136").Value = 10558.227
$wsWVR.Range("J136").Value = 856347.0600000001
$wsWVR.Range("K136").Value = 31674.681
$wsWVR.Range("L136").Value = 2569041.18
$wsWVR.Range("M136").Value = -29124.681
$wsWVR.Range("N136").Value = -2574141.18
